# Insert a new weekly price record for Achicoria (Vega Modelo de Temuco)
# at row 134, pushing the existing rows 134-158 down to 135-159.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(134).Insert()

$ws.Range("A134").Value = 10
$ws.Range("B134").Value = "Vega Modelo de Temuco"
$ws.Range("C134").Value = "La Araucanía"
$ws.Range("D134").Value = 45244
$ws.Range("E134").Value = 9
$ws.Range("F134").Value = 100112010
$ws.Range("G134").Value = "Achicoria"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 35
$ws.Range("K134").Value = 10000
$ws.Range("L134").Value = 10000
$ws.Range("M134").Value = 10000
$ws.Range("N134").Value = "$/caja 18 unidades"
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 556
$ws.Range("Q134").Value = 18
$ws.Range("R134").Value = "Hortaliza"
